$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("I3").Value = 4.5
$ws.Range("J3").Value = 2.25
$ws.Range("Q3").Value = 1.53
$ws.Range("R3").Value = 2.4
$ws.Range("W3").Value = 10
$ws.Range("X3").Value = 10
$ws.Range("AK3").Value = 41
$ws.Range("AL3").Value = 29
$ws.Range("AY3").Value = 23
$ws.Range("L4").Value = 4
$ws.Range("Q4").Value = 2.05
$ws.Range("R4").Value = 1.75
$ws.Range("S4").Value = 1.44
$ws.Range("T4").Value = 2.63
$ws.Range("U4").Value = 1.8
$ws.Range("V4").Value = 1.91
$ws.Range("W4").Value = 7.5
$ws.Range("X4").Value = 10
$ws.Range("Y4").Value = 9.5
$ws.Range("AC4").Value = 9
$ws.Range("AF4").Value = 51
$ws.Range("AG4").Value = 251
$ws.Range("AL4").Value = 29
$ws.Range("AR4").Value = 67
$ws.Range("AT4").Value = 2.63
$ws.Range("AY4").Value = 29
$ws.Range("AL5").Value = 26
$ws.Range("G12").Value = 5.75
$ws.Range("I12").Value = 1.53
$ws.Range("J12").Value = 6
$ws.Range("Q12").Value = 1.95
$ws.Range("R12").Value = 1.9
$ws.Range("W12").Value = 15
$ws.Range("Z12").Value = 67
$ws.Range("AD12").Value = 7.5
$ws.Range("AF12").Value = 51
$ws.Range("AO12").Value = 34
$ws.Range("AQ12").Value = 126
$ws.Range("AS12").Value = 500
$ws.Range("G13").Value = 2.6
$ws.Range("I13").Value = 2.6
$ws.Range("L13").Value = 3.25
$ws.Range("Q13").Value = 1.93
$ws.Range("R13").Value = 1.93
$ws.Range("AN13").Value = 4.75
$ws.Range("G14").Value = 1.3
$ws.Range("H14").Value = 5.25
$ws.Range("I14").Value = 9.5
$ws.Range("J14").Value = 1.8
$ws.Range("K14").Value = 2.5
$ws.Range("L14").Value = 8.5
$ws.Range("O14").Value = 1.2
$ws.Range("P14").Value = 4.33
$ws.Range("W14").Value = 7
$ws.Range("Y14").Value = 9
$ws.Range("Z14").Value = 8
$ws.Range("AD14").Value = 10
$ws.Range("AG14").Value = 1250
$ws.Range("AH14").Value = 23
$ws.Range("AJ14").Value = 29
$ws.Range("AK14").Value = 126
$ws.Range("AM14").Value = 67
$ws.Range("AN14").Value = 3.2
$ws.Range("AO14").Value = 6
$ws.Range("AU14").Value = 10
$ws.Range("AW14").Value = 10
$ws.Range("O15").Value = 1.36
$ws.Range("P15").Value = 3
$ws.Range("Q15").Value = 2.25
$ws.Range("R15").Value = 1.62
$ws.Range("G16").Value = 2.3
$ws.Range("I16").Value = 3.5
$ws.Range("J16").Value = 3.2
$ws.Range("K16").Value = 1.91
$ws.Range("M16").Value = 1.13
$ws.Range("N16").Value = 6
$ws.Range("O16").Value = 1.53
$ws.Range("P16").Value = 2.38
$ws.Range("Q16").Value = 2.7
$ws.Range("R16").Value = 1.44
$ws.Range("W16").Value = 5.5
$ws.Range("X16").Value = 9.5
$ws.Range("Y16").Value = 10
$ws.Range("Z16").Value = 21
$ws.Range("AC16").Value = 6
$ws.Range("AE16").Value = 21
$ws.Range("AN16").Value = 4
$ws.Range("AP16").Value = 29
$ws.Range("AR16").Value = 81
$ws.Range("AS16").Value = 301
$ws.Range("G17").Value = 1.83
$ws.Range("H17").Value = 4
$ws.Range("I17").Value = 3.5
$ws.Range("J17").Value = 2.3
$ws.Range("K17").Value = 2.47
$ws.Range("L17").Value = 3.75
$ws.Range("M17").Value = 1.03
$ws.Range("N17").Value = 9.25
$ws.Range("O17").Value = 1.16
$ws.Range("P17").Value = 4.5
$ws.Range("Q17").Value = 1.52
$ws.Range("R17").Value = 2.37
$ws.Range("S17").Value = 1.25
$ws.Range("T17").Value = 3.55
$ws.Range("U17").Value = 1.52
$ws.Range("V17").Value = 2.37
$ws.Range("W17").Value = 10.75
$ws.Range("X17").Value = 11
$ws.Range("Y17").Value = 8.5
$ws.Range("Z17").Value = 16.5
$ws.Range("AB17").Value = 18.5
$ws.Range("AC17").Value = 9.25
$ws.Range("AD17").Value = 8.25
$ws.Range("AE17").Value = 12.5
$ws.Range("AF17").Value = 40
$ws.Range("AG17").Value = 200
$ws.Range("AH17").Value = 15
$ws.Range("AI17").Value = 22
$ws.Range("AK17").Value = 50
$ws.Range("AL17").Value = 26
$ws.Range("AM17").Value = 27
$ws.Range("AN17").Value = 4.15
$ws.Range("AO17").Value = 8.75
$ws.Range("AP17").Value = 14
$ws.Range("AQ17").Value = 26
$ws.Range("AS17").Value = 120
$ws.Range("AT17").Value = 3.55
$ws.Range("AU17").Value = 6.6
$ws.Range("AV17").Value = 40
$ws.Range("AX17").Value = 17.5
$ws.Range("AY17").Value = 19.5
$ws.Range("AZ17").Value = 75
$ws.Range("BA17").Value = 90
$ws.Range("BB17").Value = 175
$ws.Range("BC17").Value = 500
